# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newer scrape counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 3021
$ws1.Range("F5").Value = 465
$ws1.Range("F6").Value = 42
$ws1.Range("F7").Value = 30
$ws1.Range("F9").Value = 9
$ws1.Range("F10").Value = 14433
$ws1.Range("F11").Value = 154
$ws1.Range("F12").Value = 119
$ws1.Range("F13").Value = 5769
$ws1.Range("F17").Value = 60
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 80
$ws1.Range("F21").Value = 178
$ws1.Range("F22").Value = 789
$ws1.Range("F23").Value = 2931
$ws1.Range("F25").Value = 10553
$ws1.Range("F29").Value = 3733

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 3021
$ws4.Range("F6").Value = 465
$ws4.Range("F7").Value = 42
$ws4.Range("F8").Value = 30
$ws4.Range("F10").Value = 9
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 154
$ws4.Range("F13").Value = 119
$ws4.Range("F14").Value = 5769
$ws4.Range("F18").Value = 60
$ws4.Range("F20").Value = 12
$ws4.Range("F21").Value = 80
$ws4.Range("F22").Value = 178
$ws4.Range("F23").Value = 789
$ws4.Range("F24").Value = 2931
$ws4.Range("F27").Value = 10553
$ws4.Range("F31").Value = 3733

$wb.Save()
